# Commit: "Added country article in db"
#
# The country name is now stored in the database together with its
# Greek grammatical article (e.g. "στην Ιταλία" instead of just
# "Ιταλία"), so the template no longer needs to hard-code the
# article/preposition ("στη χώρα: «...»") or the quoting guillemets
# around the ${country} placeholder. Two occurrences in the document
# are updated accordingly.

$d = $word.ActiveDocument

# --- Occurrence 1 ----------------------------------------------------
# Before: ' στη χώρα: «${country}»,  στο πλαίσιο  του Προγράμματος '
# After : ' ${country},  στο πλαίσιο  του  Προγράμματος '
$rng = $d.Content
$found1a = $rng.Find.Execute(' στη χώρα: «', $true, $false, $false, $false, $false, $true, 1, $false, ' ', 2)
if (-not $found1a) {
    throw "Occurrence 1a (' στη χώρα: «') not found"
}

$rng = $d.Content
$found1b = $rng.Find.Execute('}»,  στο πλαίσιο  του ', $true, $false, $false, $false, $false, $true, 1, $false, '},  στο πλαίσιο  του  ', 2)
if (-not $found1b) {
    throw "Occurrence 1b ('}»,  στο πλαίσιο  του ') not found"
}

# --- Occurrence 2 ----------------------------------------------------
# Before: 'εταιρικό σχολείο της χώρας «${country}» από '
# After : 'εταιρικό σχολείο ${country} από  '
$rng = $d.Content
$found2a = $rng.Find.Execute('της χώρας «', $true, $false, $false, $false, $false, $true, 1, $false, '', 2)
if (-not $found2a) {
    throw "Occurrence 2a ('της χώρας «') not found"
}

$rng = $d.Content
$found2b = $rng.Find.Execute('}» από ', $true, $false, $false, $false, $false, $true, 1, $false, '} από  ', 2)
if (-not $found2b) {
    throw "Occurrence 2b ('}» από ') not found"
}

Write-Output "Replacements applied: $found1a, $found1b, $found2a, $found2b"
